$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find last used row in column A (data starts at row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Update column C ("Förändrad") for every data row from 45177 to 45178
$ws.Range("C2:C$lastRow").Value = 45178
